$wb = $excel.ActiveWorkbook
$khushWiki = $wb.Worksheets.Item("Khush Wiki")
$ws = $wb.Worksheets.Add($khushWiki)
$ws.Name = "Khush Skitter"

$ws.Range("A1").Value = "Size"
$ws.Range("B1").Value = "Number of Cliques"

$sizes = 2..67
$values = @(2319807,3171609,1823321,939336,684873,598284,588889,608937,665661,728098,798073,877282,945194,980831,939987,839330,729601,639413,600192,611976,640890,673924,706753,753633,818353,892719,955212,999860,1034106,1055653,1017560,946717,878552,809485,744634,663650,583922,520239,474301,420796,367879,321829,275995,222461,158352,99522,62437,39822,30011,25637,17707,9514,3737,2042,1080,546,449,447,405,283,242,146,84,49,22,4)

for ($i = 0; $i -lt $sizes.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $sizes[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$co = $ws.ChartObjects().Add(100, 20, 400, 300)
$chart = $co.Chart
$chart.ChartType = 51
$chart.SetSourceData($ws.Range("A1:B67"))
$chart.HasTitle = $false
Write-Output "done"
